$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.955.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = "'1.858.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.23%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = "'311.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = "'0.5135"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.76%  '
$ws.Range("D8").Value = "'0.3807"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").Value = "'0.08263"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.84%  '
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").Value = "'1.107"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").Value = "'6.178"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").Value = "'20.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = "'1.858.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").Value = "'7.187"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = "'90.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").Value = "'0.06596"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").Value = "'17.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = "'5.993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.90%  '
$ws.Range("D23").Value = "'27.981.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = "'11.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.25%  '
$ws.Range("D25").Value = "'2.216"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.53%  '
$ws.Range("D26").Value = "'2.580"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").Value = "'2.074.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.24%  '
$ws.Range("D28").Value = "'156.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").Value = "'20.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.13%  '
$ws.Range("D30").Value = "'124.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").Value = "'1.036"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.45%  '
$ws.Range("D33").Value = "'5.585"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = "'3.601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").Value = "'9.505"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.88%  '
$ws.Range("D36").Value = "'0.06518"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").Value = "'0.02409"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = "'0.2174"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("D39").Value = "'1.204"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").Value = "'0.6426"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").Value = "'1.232"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.19%  '
$ws.Range("D42").Value = "'11.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.68%  '
$ws.Range("D43").Value = "'4.869"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.11%  '
$ws.Range("D44").Value = "'0.6106"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("D45").Value = "'13.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("D46").Value = "'1.276"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("D47").Value = "'3.649"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("D48").Value = "'1.970"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.86%  '
$ws.Range("D49").Value = "'1.203"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").Value = "'120.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").Value = "'79.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.95%  '
